$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new TPM-derived values, and add new rows 5-10
# reflecting the newly introduced "ECs" cluster in the Siglec1-Spn LR-pair analysis.

# Row 2
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Siglec1"
$ws.Range("C2").Value = "Spn"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.216469
$ws.Range("H2").Value = 0.6494070000000001
$ws.Range("I2").Value = 0.005147493274683087
$ws.Range("J2").Value = 0.005147493274683087
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01649366666666667
$ws.Range("N2").Value = 0.049481
$ws.Range("O2").Value = 0.005209935740510855
$ws.Range("P2").Value = 0.005209935740510855
$ws.Range("Q2").Value = 0.003570367529666667
$ws.Range("R2").Value = 0.032133307767
$ws.Range("S2").Value = [double]"2.681810918581067E-05"
$ws.Range("T2").Value = [double]"2.681810918581067E-05"

# Row 3
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Siglec1"
$ws.Range("C3").Value = "Spn"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.216469
$ws.Range("H3").Value = 0.6494070000000001
$ws.Range("I3").Value = 0.005147493274683087
$ws.Range("J3").Value = 0.005147493274683087
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.133398
$ws.Range("N3").Value = 0.400194
$ws.Range("O3").Value = 0.04213708340045676
$ws.Range("P3").Value = 0.04213708340045676
$ws.Range("Q3").Value = 0.028876531662
$ws.Range("R3").Value = 0.259888784958
$ws.Range("S3").Value = [double]"0.0002169003534186115"
$ws.Range("T3").Value = [double]"0.0002169003534186115"

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Siglec1"
$ws.Range("C4").Value = "Spn"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.216469
$ws.Range("H4").Value = 0.6494070000000001
$ws.Range("I4").Value = 0.005147493274683087
$ws.Range("J4").Value = 0.005147493274683087
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.015918333333333
$ws.Range("N4").Value = 9.047754999999999
$ws.Range("O4").Value = 0.9526529808590324
$ws.Range("P4").Value = 0.9526529808590325
$ws.Range("Q4").Value = 0.6528528256983334
$ws.Range("R4").Value = 5.875675431285
$ws.Range("S4").Value = 0.004903774812078664
$ws.Range("T4").Value = 0.004903774812078665

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Siglec1"
$ws.Range("C5").Value = "Spn"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.03281033333333333
$ws.Range("H5").Value = 0.09843099999999999
$ws.Range("I5").Value = 0.0007802085757011102
$ws.Range("J5").Value = 0.0007802085757011101
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01649366666666667
$ws.Range("N5").Value = 0.049481
$ws.Range("O5").Value = 0.005209935740510855
$ws.Range("P5").Value = 0.005209935740510855
$ws.Range("Q5").Value = 0.0005411627012222222
$ws.Range("R5").Value = 0.004870464310999999
$ws.Range("S5").Value = [double]"4.064836543598283E-06"
$ws.Range("T5").Value = [double]"4.064836543598282E-06"

# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Siglec1"
$ws.Range("C6").Value = "Spn"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.03281033333333333
$ws.Range("H6").Value = 0.09843099999999999
$ws.Range("I6").Value = 0.0007802085757011102
$ws.Range("J6").Value = 0.0007802085757011101
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.133398
$ws.Range("N6").Value = 0.400194
$ws.Range("O6").Value = 0.04213708340045676
$ws.Range("P6").Value = 0.04213708340045676
$ws.Range("Q6").Value = 0.004376832845999999
$ws.Range("R6").Value = 0.039391495614
$ws.Range("S6").Value = [double]"3.287571382406926E-05"
$ws.Range("T6").Value = [double]"3.287571382406926E-05"

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Siglec1"
$ws.Range("C7").Value = "Spn"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.03281033333333333
$ws.Range("H7").Value = 0.09843099999999999
$ws.Range("I7").Value = 0.0007802085757011102
$ws.Range("J7").Value = 0.0007802085757011101
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.015918333333333
$ws.Range("N7").Value = 9.047754999999999
$ws.Range("O7").Value = 0.9526529808590324
$ws.Range("P7").Value = 0.9526529808590325
$ws.Range("Q7").Value = 0.09895328582277776
$ws.Range("R7").Value = 0.8905795724049997
$ws.Range("S7").Value = 0.0007432680253334426
$ws.Range("T7").Value = 0.0007432680253334426

# Row 8
$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Siglec1"
$ws.Range("C8").Value = "Spn"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 41.80400533333333
$ws.Range("H8").Value = 125.412016
$ws.Range("I8").Value = 0.9940722981496158
$ws.Range("J8").Value = 0.9940722981496158
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01649366666666667
$ws.Range("N8").Value = 0.049481
$ws.Range("O8").Value = 0.005209935740510855
$ws.Range("P8").Value = 0.005209935740510855
$ws.Range("Q8").Value = 0.6895013292995555
$ws.Range("R8").Value = 6.205511963695999
$ws.Range("S8").Value = 0.005179052794781445
$ws.Range("T8").Value = 0.005179052794781445

# Row 9
$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Siglec1"
$ws.Range("C9").Value = "Spn"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 41.80400533333333
$ws.Range("H9").Value = 125.412016
$ws.Range("I9").Value = 0.9940722981496158
$ws.Range("J9").Value = 0.9940722981496158
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.133398
$ws.Range("N9").Value = 0.400194
$ws.Range("O9").Value = 0.04213708340045676
$ws.Range("P9").Value = 0.04213708340045676
$ws.Range("Q9").Value = 5.576570703455999
$ws.Range("R9").Value = 50.189136331104
$ws.Range("S9").Value = 0.04188730733321408
$ws.Range("T9").Value = 0.04188730733321409

# Row 10
$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Siglec1"
$ws.Range("C10").Value = "Spn"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 41.80400533333333
$ws.Range("H10").Value = 125.412016
$ws.Range("I10").Value = 0.9940722981496158
$ws.Range("J10").Value = 0.9940722981496158
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.015918333333333
$ws.Range("N10").Value = 9.047754999999999
$ws.Range("O10").Value = 0.9526529808590324
$ws.Range("P10").Value = 0.9526529808590325
$ws.Range("Q10").Value = 126.0774660915644
$ws.Range("R10").Value = 1134.69719482408
$ws.Range("S10").Value = 0.9470059380216204
$ws.Range("T10").Value = 0.9470059380216205
